$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.389
$ws.Range("C10").Value = -13.09
$ws.Range("C12").Value = -11.066
$ws.Range("E15").Value = 16.353
$ws.Range("C18").Value = -11.681
$ws.Range("E20").Value = 16.355
$ws.Range("E29").Value = 17.05
$ws.Range("E30").Value = 16.436
$ws.Range("E31").Value = 16.225
$ws.Range("C37").Value = -13.332
$ws.Range("E40").Value = 16.789
$ws.Range("C55").Value = -13.916
$ws.Range("C68").Value = -11.262
$ws.Range("E68").Value = 17.006
$ws.Range("E76").Value = 16.697
$ws.Range("C77").Value = -12.824
$ws.Range("C78").Value = -12.798
$ws.Range("E87").Value = 16.317
$ws.Range("E88").Value = 16.213
$ws.Range("E96").Value = 16.325
$ws.Range("E98").Value = 16.299
$ws.Range("E101").Value = 16.625
$ws.Range("E102").Value = 16.64
